$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one more day of data as the new last row (row 64), right after the
# existing last row (63). The new row's B:J values are identical to the
# previous row, so copy the whole row (values + formatting, e.g. column A's
# date style) down and then fix up the date in column A.
$lastRow = 63
$newRow = $lastRow + 1

$srcRange = $ws.Range("A" + $lastRow + ":J" + $lastRow)
$dstRange = $ws.Range("A" + $newRow + ":J" + $newRow)
$srcRange.Copy($dstRange)

$ws.Cells.Item($newRow, 1).Value = 45620
